$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend header row (row 1) with two new columns P (16) and Q (17) ---
# Copy the formatting of O1 (bold, bordered, centered header style) onto P1:Q1
$ws.Cells.Item(1, 15).Copy()
$ws.Range($ws.Cells.Item(1, 16), $ws.Cells.Item(1, 17)).PasteSpecial(-4122)
$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15

# --- Update data rows 2-25 ---
# For every row: swap values between columns I/K and M/O, then add new
# columns P and Q (both value 2), matching the parallel-line contingency table.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> was 1
    $ws.Cells.Item($r, 11).Value = 1   # K -> was 2
    $ws.Cells.Item($r, 13).Value = 2   # M -> was 1
    $ws.Cells.Item($r, 15).Value = 1   # O -> was 2
    $ws.Cells.Item($r, 16).Value = 2   # P -> new
    $ws.Cells.Item($r, 17).Value = 2   # Q -> new
}
